$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 18183562
$ws.Cells.Item(100, 9).Value = 980
$ws.Cells.Item(100, 10).Value = 40002660
$ws.Cells.Item(100, 11).Value = 980
$ws.Cells.Item(100, 12).Value = 40002660
$ws.Cells.Item(100, 13).Value = -439
$ws.Cells.Item(100, 14).Value = -40003742
$ws.Cells.Item(113, 8).Value = 6668688.5
$ws.Cells.Item(113, 9).Value = 9092727
$ws.Cells.Item(113, 10).Value = 2583.25
$ws.Cells.Item(113, 11).Value = 9092727
$ws.Cells.Item(113, 12).Value = 2583.25
$ws.Cells.Item(113, 13).Value = -9089473
$ws.Cells.Item(113, 14).Value = -9091.25
$ws.Cells.Item(116, 8).Value = 10656624
$ws.Cells.Item(116, 9).Value = 6668646.5
$ws.Cells.Item(116, 10).Value = 13979938
$ws.Cells.Item(116, 11).Value = 6668646.5
$ws.Cells.Item(116, 12).Value = 13979938
$ws.Cells.Item(116, 13).Value = -6665204.5
$ws.Cells.Item(116, 14).Value = -13986822
$ws.Cells.Item(137, 8).Value = 11583985
$ws.Cells.Item(137, 9).Value = 806.5833
$ws.Cells.Item(137, 10).Value = 104249416
$ws.Cells.Item(137, 11).Value = 2419.7499
$ws.Cells.Item(137, 12).Value = 312748248
$ws.Cells.Item(137, 13).Value = 130.2501000000002
$ws.Cells.Item(137, 14).Value = -312753348
$ws.Cells.Item(138, 8).Value = 2734.7214
$ws.Cells.Item(138, 9).Value = 2079.5757
$ws.Cells.Item(138, 10).Value = 3204.7173
$ws.Cells.Item(138, 11).Value = 6238.7271
$ws.Cells.Item(138, 12).Value = 9614.151899999999
$ws.Cells.Item(138, 13).Value = -1098.7271
$ws.Cells.Item(138, 14).Value = -19894.1519
$ws.Cells.Item(141, 8).Value = 1180.1459
$ws.Cells.Item(141, 9).Value = 910.7619
$ws.Cells.Item(141, 10).Value = 3065.8333
$ws.Cells.Item(141, 11).Value = 2732.2857
$ws.Cells.Item(141, 12).Value = 9197.499899999999
$ws.Cells.Item(141, 13).Value = 2447.7143
$ws.Cells.Item(141, 14).Value = -19557.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 313546.16
$ws.Cells.Item(45, 9).Value = 345837.12
$ws.Cells.Item(45, 10).Value = 1400
$ws.Cells.Item(45, 11).Value = 345837.12
$ws.Cells.Item(45, 12).Value = 1400
$ws.Cells.Item(45, 13).Value = -345460.12
$ws.Cells.Item(45, 14).Value = -2154
$ws.Cells.Item(74, 8).Value = 19608844
$ws.Cells.Item(74, 9).Value = 16667169
$ws.Cells.Item(74, 10).Value = 41671400
$ws.Cells.Item(74, 11).Value = 16667169
$ws.Cells.Item(74, 12).Value = 41671400
$ws.Cells.Item(74, 13).Value = -16666295
$ws.Cells.Item(74, 14).Value = -41673148
$ws.Cells.Item(77, 8).Value = 19608844
$ws.Cells.Item(77, 9).Value = 16667169
$ws.Cells.Item(77, 10).Value = 41671400
$ws.Cells.Item(77, 11).Value = 83335845
$ws.Cells.Item(77, 12).Value = 208357000
$ws.Cells.Item(77, 13).Value = -83331477
$ws.Cells.Item(77, 14).Value = -208365736
$ws.Cells.Item(110, 8).Value = 653.5925999999999
$ws.Cells.Item(110, 9).Value = 566.2727
$ws.Cells.Item(110, 10).Value = 1037.8
$ws.Cells.Item(110, 11).Value = 566.2727
$ws.Cells.Item(110, 12).Value = 1037.8
$ws.Cells.Item(110, 13).Value = 1478.7273
$ws.Cells.Item(110, 14).Value = -5127.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 263496.03
$ws.Cells.Item(107, 9).Value = 435060.4
$ws.Cells.Item(107, 10).Value = 430.66666
$ws.Cells.Item(107, 11).Value = 435060.4
$ws.Cells.Item(107, 12).Value = 430.66666
$ws.Cells.Item(107, 13).Value = -433140.4
$ws.Cells.Item(107, 14).Value = -4270.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2323.6365
$ws.Cells.Item(16, 9).Value = 2032.5
$ws.Cells.Item(16, 10).Value = 3100
$ws.Cells.Item(16, 11).Value = 2032.5
$ws.Cells.Item(16, 12).Value = 3100
$ws.Cells.Item(16, 13).Value = -1745.5
$ws.Cells.Item(16, 14).Value = -3674
$ws.Cells.Item(31, 8).Value = 1897548.8
$ws.Cells.Item(31, 9).Value = 1104.0435
$ws.Cells.Item(31, 10).Value = 6259371.5
$ws.Cells.Item(31, 11).Value = 1104.0435
$ws.Cells.Item(31, 12).Value = 6259371.5
$ws.Cells.Item(31, 13).Value = -809.0435
$ws.Cells.Item(31, 14).Value = -6259961.5
$ws.Cells.Item(34, 8).Value = 1897548.8
$ws.Cells.Item(34, 9).Value = 1104.0435
$ws.Cells.Item(34, 10).Value = 6259371.5
$ws.Cells.Item(34, 11).Value = 1104.0435
$ws.Cells.Item(34, 12).Value = 6259371.5
$ws.Cells.Item(34, 13).Value = -902.0435
$ws.Cells.Item(34, 14).Value = -6259775.5
$ws.Cells.Item(107, 8).Value = 643.6842
$ws.Cells.Item(107, 9).Value = 235.2
$ws.Cells.Item(107, 10).Value = 789.5714
$ws.Cells.Item(107, 11).Value = 235.2
$ws.Cells.Item(107, 12).Value = 789.5714
$ws.Cells.Item(107, 13).Value = 1684.8
$ws.Cells.Item(107, 14).Value = -4629.5714
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws.Cells.Item(111, 8).Value = 42000
$ws.Cells.Item(111, 10).Value = 42000
$ws.Cells.Item(111, 12).Value = 42000
$ws.Cells.Item(111, 14).Value = -50180
$ws.Cells.Item(112, 8).Value = 34799
$ws.Cells.Item(112, 10).Value = 34799
$ws.Cells.Item(112, 12).Value = 34799
$ws.Cells.Item(112, 14).Value = -37753
$ws.Cells.Item(113, 8).Value = 2323.6365
$ws.Cells.Item(113, 9).Value = 2032.5
$ws.Cells.Item(113, 10).Value = 3100
$ws.Cells.Item(113, 11).Value = 2032.5
$ws.Cells.Item(113, 12).Value = 3100
$ws.Cells.Item(113, 13).Value = 137.5
$ws.Cells.Item(113, 14).Value = -7440
$ws.Cells.Item(114, 8).Value = 20684
$ws.Cells.Item(114, 10).Value = 20684
$ws.Cells.Item(114, 12).Value = 20684
$ws.Cells.Item(114, 14).Value = -29362
$ws.Cells.Item(115, 8).Value = 500005000
$ws.Cells.Item(115, 9).Value = 10000
$ws.Cells.Item(115, 10).Value = 1000000000
$ws.Cells.Item(115, 11).Value = 10000
$ws.Cells.Item(115, 12).Value = 1000000000
$ws.Cells.Item(115, 13).Value = -8825
$ws.Cells.Item(115, 14).Value = -1000002350
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()
$ws.Cells.Item(117, 8).Value = 15000
$ws.Cells.Item(117, 10).Value = 15000
$ws.Cells.Item(117, 12).Value = 15000
$ws.Cells.Item(117, 14).Value = -24178
$ws.Cells.Item(118, 8).Value = 46000
$ws.Cells.Item(118, 10).Value = 46000
$ws.Cells.Item(118, 12).Value = 46000
$ws.Cells.Item(118, 14).Value = -49314
$ws.Cells.Item(119, 8).Value = 28000
$ws.Cells.Item(119, 10).Value = 28000
$ws.Cells.Item(119, 12).Value = 28000
$ws.Cells.Item(119, 14).Value = -37676

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 950426.4399999999
$ws.Cells.Item(107, 9).Value = 2564221
$ws.Cells.Item(107, 10).Value = 1135.5294
$ws.Cells.Item(107, 11).Value = 7692663
$ws.Cells.Item(107, 12).Value = 3406.5882
$ws.Cells.Item(107, 13).Value = -7690743
$ws.Cells.Item(107, 14).Value = -7246.5882
$ws.Cells.Item(108, 8).Value = 814.5714
$ws.Cells.Item(108, 9).Value = 355
$ws.Cells.Item(108, 10).Value = 998.4
$ws.Cells.Item(108, 11).Value = 1065
$ws.Cells.Item(108, 12).Value = 2995.2
$ws.Cells.Item(108, 13).Value = 1815
$ws.Cells.Item(108, 14).Value = -8755.200000000001
$ws.Cells.Item(109, 8).Value = 3669.0466
$ws.Cells.Item(109, 9).Value = 500
$ws.Cells.Item(109, 10).Value = 3823.634
$ws.Cells.Item(109, 11).Value = 1500
$ws.Cells.Item(109, 12).Value = 11470.902
$ws.Cells.Item(109, 13).Value = -460
$ws.Cells.Item(109, 14).Value = -13550.902
$ws.Cells.Item(111, 8).Value = 125001020
$ws.Cells.Item(111, 9).Value = 250000140
$ws.Cells.Item(111, 10).Value = 1900
$ws.Cells.Item(111, 11).Value = 750000420
$ws.Cells.Item(111, 12).Value = 5700
$ws.Cells.Item(111, 13).Value = -749997353
$ws.Cells.Item(111, 14).Value = -11834
$ws.Cells.Item(112, 8).Value = 3106
$ws.Cells.Item(112, 10).Value = 3106
$ws.Cells.Item(112, 12).Value = 9318
$ws.Cells.Item(112, 14).Value = -11534
$ws.Cells.Item(113, 8).Value = 937293.7
$ws.Cells.Item(113, 9).Value = 668
$ws.Cells.Item(113, 10).Value = 2033831
$ws.Cells.Item(113, 11).Value = 2004
$ws.Cells.Item(113, 12).Value = 6101493
$ws.Cells.Item(113, 13).Value = 166
$ws.Cells.Item(113, 14).Value = -6105833
$ws.Cells.Item(114, 8).Value = 980.2222
$ws.Cells.Item(114, 9).Value = 301.83334
$ws.Cells.Item(114, 10).Value = 1522.9333
$ws.Cells.Item(114, 11).Value = 905.5000200000001
$ws.Cells.Item(114, 12).Value = 4568.7999
$ws.Cells.Item(114, 13).Value = 2348.49998
$ws.Cells.Item(114, 14).Value = -11076.7999
$ws.Cells.Item(115, 8).Value = 3755
$ws.Cells.Item(115, 9).Value = 2000
$ws.Cells.Item(115, 10).Value = 3914.5454
$ws.Cells.Item(115, 11).Value = 6000
$ws.Cells.Item(115, 12).Value = 11743.6362
$ws.Cells.Item(115, 13).Value = -4825
$ws.Cells.Item(115, 14).Value = -14093.6362
$ws.Cells.Item(116, 8).Value = 1685216.6
$ws.Cells.Item(116, 9).Value = 2000
$ws.Cells.Item(116, 10).Value = 2526825
$ws.Cells.Item(116, 11).Value = 6000
$ws.Cells.Item(116, 12).Value = 7580475
$ws.Cells.Item(116, 13).Value = -2558
$ws.Cells.Item(116, 14).Value = -7587359
$ws.Cells.Item(117, 8).Value = 1649
$ws.Cells.Item(117, 9).Value = 586
$ws.Cells.Item(117, 10).Value = 1836.5883
$ws.Cells.Item(117, 11).Value = 1758
$ws.Cells.Item(117, 12).Value = 5509.7649
$ws.Cells.Item(117, 13).Value = 1684
$ws.Cells.Item(117, 14).Value = -12393.7649
$ws.Cells.Item(118, 8).Value = 5000
$ws.Cells.Item(118, 9).Value = 5000
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 15000
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = -13757
$ws.Cells.Item(118, 14).ClearContents()
$ws.Cells.Item(119, 8).Value = 4099.2
$ws.Cells.Item(119, 9).Value = 1000
$ws.Cells.Item(119, 10).Value = 4874
$ws.Cells.Item(119, 11).Value = 3000
$ws.Cells.Item(119, 12).Value = 14622
$ws.Cells.Item(119, 13).Value = 1838
$ws.Cells.Item(119, 14).Value = -24298
$ws.Cells.Item(120, 8).Value = 41669732
$ws.Cells.Item(120, 10).Value = 3500
$ws.Cells.Item(120, 12).Value = 10500
$ws.Cells.Item(120, 14).Value = -20176

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 119.875
$ws.Cells.Item(107, 9).Value = 119.875
$ws.Cells.Item(107, 11).Value = 119.875
$ws.Cells.Item(107, 13).Value = 1800.125
$ws.Cells.Item(113, 8).Value = 84333.336
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 84333.336
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 84333.336
$ws.Cells.Item(113, 14).Value = -88673.336
$ws.Cells.Item(113, 13).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(107, 8).Value = 6510
$ws.Cells.Item(107, 9).Value = 6510
$ws.Cells.Item(107, 11).Value = 6510
$ws.Cells.Item(107, 13).Value = -4590
$ws.Cells.Item(136, 8).Value = 2850069.5
$ws.Cells.Item(136, 9).Value = 3004096.2
$ws.Cells.Item(136, 10).Value = 572.5
$ws.Cells.Item(136, 11).Value = 9012288.600000001
$ws.Cells.Item(136, 12).Value = 1717.5
$ws.Cells.Item(136, 13).Value = -9009738.600000001
$ws.Cells.Item(136, 14).Value = -6817.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 6375.353
$ws.Cells.Item(100, 9).Value = 8790.833000000001
$ws.Cells.Item(100, 10).Value = 578.2
$ws.Cells.Item(100, 11).Value = 17581.666
$ws.Cells.Item(100, 12).Value = 1156.4
$ws.Cells.Item(100, 13).Value = -17040.666
$ws.Cells.Item(100, 14).Value = -2238.4
$ws.Cells.Item(136, 8).Value = 1280.3334
$ws.Cells.Item(136, 9).Value = 475.5
$ws.Cells.Item(136, 10).Value = 2158.3333
$ws.Cells.Item(136, 11).Value = 1426.5
$ws.Cells.Item(136, 12).Value = 6474.999899999999
$ws.Cells.Item(136, 13).Value = 1123.5
$ws.Cells.Item(136, 14).Value = -11574.9999
